# Append a new submission row (row 12) to the log sheet, mirroring the
# shape of the existing rows (same person/camp/trip/vehicle/org as row 6,
# just a later timestamp).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A is an empty "note" field. A plain "" assignment clears the cell
# entirely instead of leaving a blank-but-present text cell, so use a
# lone quote prefix to force an empty *text* cell (matches A2/A3/A6-A11).
$ws.Range("A12").Value = "'"

$ws.Range("B12").Value = "أحمد شريم"

# "2323" / "C2" look numeric/alphanumeric - a plain string assignment of a
# pure-digit string gets auto-coerced to a Number by the engine, so force
# text storage with a quote prefix (same trick real users use in Excel).
$ws.Range("C12").Value = "'2323"

$ws.Range("D12").Value = "ايتا"
$ws.Range("E12").Value = "الرحلة 2"
$ws.Range("F12").Value = "'C2"
$ws.Range("G12").Value = "NRC"
$ws.Range("H12").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٥٧:٢٥ م"
